$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value = 1
$ws.Range("F21").Value = -2
$ws.Range("F27").Value = 0
$ws.Range("F33").Value = 3
$ws.Range("F34").Value = 2
$ws.Range("F36").Value = 3
$ws.Range("F50").Value = -1
$ws.Range("F54").Value = 0
$ws.Range("F63").Value = 3
$ws.Range("F64").Value = 2
$ws.Range("F67").Value = -2
$ws.Range("F73").Value = 2
